$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 304.42856
$ws.Range("I53").Value = 232.25
$ws.Range("J53").Value = 400.66666
$ws.Range("K53").Value = 232.25
$ws.Range("L53").Value = 400.66666
$ws.Range("M53").Value = 404.75
$ws.Range("N53").Value = -1674.66666
$ws.Range("H112").Value = 1790.4546
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 1869.5
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 5608.5
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -7824.5
$ws.Range("H137").Value = 5320710
$ws.Range("I137").Value = 9616591
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 28849773
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -28847223
$ws.Range("N137").Value = -11100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1000
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("H32").Value = 16954986
$ws.Range("I32").Value = 25643304
$ws.Range("K32").Value = 25643304
$ws.Range("M32").Value = -25643017
$ws.Range("H61").Value = 1958.68
$ws.Range("I61").Value = 1911.7778
$ws.Range("J61").Value = 2079.2856
$ws.Range("K61").Value = 1911.7778
$ws.Range("L61").Value = 2079.2856
$ws.Range("M61").Value = -1699.7778
$ws.Range("N61").Value = -2503.2856
$ws.Range("H62").Value = 9999
$ws.Range("J62").Value = 9999
$ws.Range("L62").Value = 9999
$ws.Range("N62").Value = -11247
$ws.Range("H63").Value = 2849.8333
$ws.Range("I63").Value = 3099.75
$ws.Range("J63").Value = 2350
$ws.Range("K63").Value = 3099.75
$ws.Range("L63").Value = 2350
$ws.Range("M63").Value = -2413.75
$ws.Range("N63").Value = -3722
$ws.Range("H65").Value = 9999
$ws.Range("J65").Value = 9999
$ws.Range("L65").Value = 29997
$ws.Range("N65").Value = -36237
$ws.Range("H66").Value = 2849.8333
$ws.Range("I66").Value = 3099.75
$ws.Range("J66").Value = 2350
$ws.Range("K66").Value = 15498.75
$ws.Range("L66").Value = 11750
$ws.Range("M66").Value = -12066.75
$ws.Range("N66").Value = -18614
$ws.Range("H74").Value = 4549.4644
$ws.Range("I74").Value = 7314.2666
$ws.Range("J74").Value = 1359.3077
$ws.Range("K74").Value = 7314.2666
$ws.Range("L74").Value = 1359.3077
$ws.Range("M74").Value = -6440.2666
$ws.Range("N74").Value = -3107.3077
$ws.Range("H75").Value = 37500
$ws.Range("J75").Value = 37500
$ws.Range("L75").Value = 37500
$ws.Range("N75").Value = -39248
$ws.Range("H77").Value = 4549.4644
$ws.Range("I77").Value = 7314.2666
$ws.Range("J77").Value = 1359.3077
$ws.Range("K77").Value = 36571.333
$ws.Range("L77").Value = 6796.538500000001
$ws.Range("M77").Value = -32203.333
$ws.Range("N77").Value = -15532.5385
$ws.Range("H78").Value = 37500
$ws.Range("J78").Value = 37500
$ws.Range("L78").Value = 112500
$ws.Range("N78").Value = -121236
$ws.Range("H81").Value = 50181
$ws.Range("J81").Value = 50181
$ws.Range("L81").Value = 50181
$ws.Range("N81").Value = -52177
$ws.Range("H84").Value = 50181
$ws.Range("J84").Value = 50181
$ws.Range("L84").Value = 150543
$ws.Range("N84").Value = -160527
$ws.Range("H110").Value = 1155.5454
$ws.Range("I110").Value = 1000.0625
$ws.Range("J110").Value = 1570.1666
$ws.Range("K110").Value = 1000.0625
$ws.Range("L110").Value = 1570.1666
$ws.Range("M110").Value = 1044.9375
$ws.Range("N110").Value = -5660.1666
$ws.Range("H122").Value = 1957.3334
$ws.Range("I122").Value = 1957.3334
$ws.Range("K122").Value = 5872.0002
$ws.Range("M122").Value = -3422.0002
$ws.Range("H132").Value = 2685.6924
$ws.Range("I132").Value = 2713.12
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 8139.36
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -5609.36
$ws.Range("N132").Value = -11060
$ws.Range("H136").Value = 1958.68
$ws.Range("I136").Value = 1911.7778
$ws.Range("J136").Value = 2079.2856
$ws.Range("K136").Value = 5735.3334
$ws.Range("L136").Value = 6237.8568
$ws.Range("M136").Value = -3185.3334
$ws.Range("N136").Value = -11337.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 3782.5
$ws.Range("I113").Value = 3782.5
$ws.Range("K113").Value = 3782.5
$ws.Range("M113").Value = -1612.5
$ws.Range("H134").Value = 2622.1667
$ws.Range("I134").Value = 2587.182
$ws.Range("K134").Value = 7761.545999999999
$ws.Range("M134").Value = -5226.545999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1408.8422
$ws.Range("I31").Value = 1436.75
$ws.Range("K31").Value = 1436.75
$ws.Range("M31").Value = -1141.75
$ws.Range("H34").Value = 1408.8422
$ws.Range("I34").Value = 1436.75
$ws.Range("K34").Value = 1436.75
$ws.Range("M34").Value = -1234.75
$ws.Range("H132").Value = 50001920
$ws.Range("I132").Value = 66668060
$ws.Range("J132").Value = 3499.6
$ws.Range("K132").Value = 200004180
$ws.Range("L132").Value = 10498.8
$ws.Range("M132").Value = -200001650
$ws.Range("N132").Value = -15558.8
$ws.Range("H134").Value = 3433.9688
$ws.Range("I134").Value = 1882.1666
$ws.Range("J134").Value = 5429.143
$ws.Range("K134").Value = 5646.4998
$ws.Range("L134").Value = 16287.429
$ws.Range("M134").Value = -3111.4998
$ws.Range("N134").Value = -21357.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 500034900
$ws.Range("J69").Value = 69800
$ws.Range("L69").Value = 69800
$ws.Range("N69").Value = -71298
$ws.Range("H72").Value = 500034900
$ws.Range("J72").Value = 69800
$ws.Range("L72").Value = 209400
$ws.Range("N72").Value = -216888
$ws.Range("H87").Value = 30354
$ws.Range("J87").Value = 30354
$ws.Range("L87").Value = 30354
$ws.Range("N87").Value = -32850
$ws.Range("H90").Value = 30354
$ws.Range("J90").Value = 30354
$ws.Range("L90").Value = 91062
$ws.Range("N90").Value = -103542
$ws.Range("H92").Value = 18900
$ws.Range("J92").Value = 18900
$ws.Range("L92").Value = 18900
$ws.Range("N92").Value = -22644
$ws.Range("H132").Value = 4360.5654
$ws.Range("I132").Value = 4618.4375
$ws.Range("J132").Value = 3771.1428
$ws.Range("K132").Value = 13855.3125
$ws.Range("L132").Value = 11313.4284
$ws.Range("M132").Value = -11325.3125
$ws.Range("N132").Value = -16373.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4521.2104
$ws.Range("I40").Value = 3928.7856
$ws.Range("J40").Value = 6180
$ws.Range("K40").Value = 3928.7856
$ws.Range("L40").Value = 6180
$ws.Range("M40").Value = -3792.7856
$ws.Range("N40").Value = -6452
$ws.Range("H132").Value = 3367.5833
$ws.Range("I132").Value = 2991.6191
$ws.Range("J132").Value = 5999.3335
$ws.Range("K132").Value = 8974.8573
$ws.Range("L132").Value = 17998.0005
$ws.Range("M132").Value = -6444.8573
$ws.Range("N132").Value = -23058.0005
$ws.Range("H136").Value = 2033.3334
$ws.Range("I136").Value = 1044.4445
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 3133.3335
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -583.3335000000002
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2264.325
$ws.Range("I122").Value = 1689
$ws.Range("J122").Value = 3042.7058
$ws.Range("K122").Value = 5067
$ws.Range("L122").Value = 9128.117400000001
$ws.Range("M122").Value = -2617
$ws.Range("N122").Value = -14028.1174
$ws.Range("H132").Value = 1898.4062
$ws.Range("I132").Value = 1455.3572
$ws.Range("J132").Value = 4999.75
$ws.Range("K132").Value = 4366.071599999999
$ws.Range("L132").Value = 14999.25
$ws.Range("M132").Value = -1836.071599999999
$ws.Range("N132").Value = -20059.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N30").ClearContents()
